$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COCKTAIL & BAR")
$ws.Activate()
Write-Host $ws.Name
